$wb = $excel.ActiveWorkbook

# --- Sheet "Sheet": add two new client rows (13 and 14) ---
$ws1 = $wb.Worksheets.Item("Sheet")

$ws1.Cells.Item(13, 1).Value = "Santiago Arango"
$ws1.Cells.Item(13, 2).Value = 12
$ws1.Cells.Item(13, 3).Value = "santi@gmail.com"
$ws1.Cells.Item(13, 4).Value = "Medellín"
$ws1.Cells.Item(13, 5).Value = 7155934

$ws1.Cells.Item(14, 1).Value = "Bayron Valdés "
$ws1.Cells.Item(14, 2).Value = 13
$ws1.Cells.Item(14, 3).Value = "bayron2813@gmail.com"
$ws1.Cells.Item(14, 4).Value = "Medellín"
$ws1.Cells.Item(14, 5).Value = 3015168866

# --- Sheet "Pedidos": add five new order rows (8 through 12) ---
$ws2 = $wb.Worksheets.Item("Pedidos")

$ws2.Cells.Item(8, 1).Value = 12
$ws2.Cells.Item(8, 2).Value = 7
$ws2.Cells.Item(8, 3).Value = "El tamaño del Diamante es 0.3 `nEl grabado del Diamante es True `nEl origen del diamante es cabello `nEl tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl corte del diamante es corazon `n"
$ws2.Cells.Item(8, 4).Value = "31/12/2022"

$ws2.Cells.Item(9, 1).Value = 13
$ws2.Cells.Item(9, 2).Value = 8
$ws2.Cells.Item(9, 3).Value = "El tamaño del Diamante es 0.8 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl corte del diamante es Corazon `n"
$ws2.Cells.Item(9, 4).Value = "23/04/2021"

$ws2.Cells.Item(10, 1).Value = 13
$ws2.Cells.Item(10, 2).Value = 9
$ws2.Cells.Item(10, 3).Value = "El tamaño del Diamante es 0.6 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `nEl corte del diamante es corazon `n"
$ws2.Cells.Item(10, 4).Value = "20/04/2022"

$ws2.Cells.Item(11, 1).Value = 10
$ws2.Cells.Item(11, 2).Value = 10
$ws2.Cells.Item(11, 3).Value = "El tamaño del Diamante es 0.5 `nEl grabado del Diamante es True `nEl origen del diamante es cenizas `n"
$ws2.Cells.Item(11, 4).Value = "20/04/2022"

$ws2.Cells.Item(12, 1).Value = 10
$ws2.Cells.Item(12, 2).Value = 11
$ws2.Cells.Item(12, 3).Value = "El tamaño del Diamante es 0.6 `nEl grabado del Diamante es True `nEl origen del diamante es Cabello `n"
$ws2.Cells.Item(12, 4).Value = "21/04/2022"
